# Cotações atualizadas - 2025-12-07
# Append a new row (93) with the latest fund quotes, mirroring the
# formatting of the preceding rows (row 92 = 2025-12-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 93

# Column A: date value, stored as a real number (days since 1899-12-30)
# with the same style as the other date cells in column A.
$ws.Cells.Item($row, 1).Value = 45998
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Columns B-E: quote values, stored as inline text (comma decimal
# separator), matching the existing rows.
$ws.Cells.Item($row, 2).Value = "22,2192"
$ws.Cells.Item($row, 3).Value = "15,9549"
$ws.Cells.Item($row, 4).Value = "15,6518"
$ws.Cells.Item($row, 5).Value = "15,6518"
